$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) - unchanged text, rewritten for determinism
$ws.Range("A1").Value = "enfermedad"
$ws.Range("B1").Value = "prevalencia"
$ws.Range("C1").Value = "year"
$ws.Range("D1").Value = "ref"
$ws.Range("E1").Value = "ref_link"

# Row 2: Diabetes Mellitus (DM)
$ws.Range("A2").Value = "DM"
$ws.Range("C2").Value = 2022
$ws.Range("D2").Value = "Fondo Colombiano de Enfermedades de Alto Costo. Situación de la enfermedad renal crónica, la hipertensión arterial y la diabetes mellitus en Colombia 2022. Colombia: CAC; 2023. Pag. 54"

# Row 3: Enfermedad Renal Crónica (ERC)
$ws.Range("A3").Value = "ERC"
$ws.Range("C3").Value = 2022
$ws.Range("D3").Value = "Fondo Colombiano de Enfermedades de Alto Costo. Situación de la enfermedad renal crónica, la hipertensión arterial y la diabetes mellitus en Colombia 2022. Colombia: CAC; 2023. Pag. 67"

# Hyperlinks first, so the generated "Hyperlink" style lands at cellXfs index 2
$url = "https://cuentadealtocosto.org/wp-content/uploads/2023/10/final-libro-erc-2022-2.pdf"
$ws.Hyperlinks.Add($ws.Range("E2"), $url, "", "", $url)
$ws.Hyperlinks.Add($ws.Range("E3"), $url, "", "", $url)

# prevalencia values must stay text (not auto-converted to numbers)
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "0.0312"
$ws.Range("B2").Style = "Normal"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "0.0154"
$ws.Range("B3").Style = "Normal"

$ws.Range("B3").Select() | Out-Null
